$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Update price (D) and volume (E) columns for existing rows ---
Set-TextValue $ws.Range("D2") "57.602.16"
Set-TextValue $ws.Range("E2") "  -2.60%  "

Set-TextValue $ws.Range("D3") "3.079.47"
Set-TextValue $ws.Range("E3") "  -1.26%  "

Set-TextValue $ws.Range("E4") "  -0.03%  "

Set-TextValue $ws.Range("D5") "523.97"
Set-TextValue $ws.Range("E5") "  -0.13%  "

Set-TextValue $ws.Range("D6") "140.41"
Set-TextValue $ws.Range("E6") "  -3.24%  "

Set-TextValue $ws.Range("E7") "  +0.28%  "

Set-TextValue $ws.Range("D8") "3.078.33"
Set-TextValue $ws.Range("E8") "  -1.22%  "

Set-TextValue $ws.Range("D9") "0.440"
Set-TextValue $ws.Range("E9") "  +0.04%  "

Set-TextValue $ws.Range("D10") "7.14"
Set-TextValue $ws.Range("E10") "  -3.73%  "

Set-TextValue $ws.Range("E11") "  -1.76%  "

Set-TextValue $ws.Range("E12") "  +1.44%  "

Set-TextValue $ws.Range("D13") "3.610.69"
Set-TextValue $ws.Range("E13") "  -1.13%  "

Set-TextValue $ws.Range("E14") "  +1.77%  "

Set-TextValue $ws.Range("E15") "  -7.12%  "

Set-TextValue $ws.Range("E16") "  -2.65%  "

Set-TextValue $ws.Range("D17") "57.656.32"
Set-TextValue $ws.Range("E17") "  -2.42%  "

Set-TextValue $ws.Range("D18") "3.079.91"
Set-TextValue $ws.Range("E18") "  -0.38%  "

Set-TextValue $ws.Range("D19") "6.06"
Set-TextValue $ws.Range("E19") "  -2.91%  "

Set-TextValue $ws.Range("E20") "  -3.54%  "

Set-TextValue $ws.Range("D21") "7.92"
Set-TextValue $ws.Range("E21") "  -4.44%  "

Set-TextValue $ws.Range("D22") "340.03"
Set-TextValue $ws.Range("E22") "  -0.37%  "

Set-TextValue $ws.Range("D23") "0.999"
Set-TextValue $ws.Range("E23") "  -0.18%  "

Set-TextValue $ws.Range("D24") "0.509"
Set-TextValue $ws.Range("E24") "  -0.48%  "

Set-TextValue $ws.Range("D25") "67.15"
Set-TextValue $ws.Range("E25") "  +1.34%  "

Set-TextValue $ws.Range("E26") "  -2.60%  "

Set-TextValue $ws.Range("E27") "  +0.18%  "

Set-TextValue $ws.Range("E28") "  -1.73%  "

Set-TextValue $ws.Range("D29") "0.999"
Set-TextValue $ws.Range("E29") "  +0.07%  "

Set-TextValue $ws.Range("D30") "6.35"
Set-TextValue $ws.Range("E30") "  -4.67%  "

Set-TextValue $ws.Range("D31") "7.19"
Set-TextValue $ws.Range("E31") "  -2.24%  "

Set-TextValue $ws.Range("E32") "  +0.70%  "

Set-TextValue $ws.Range("D33") "20.86"
Set-TextValue $ws.Range("E33") "  -1.51%  "

Set-TextValue $ws.Range("E34") "  -4.73%  "

Set-TextValue $ws.Range("D35") "158.76"
Set-TextValue $ws.Range("E35") "  +2.08%  "

Set-TextValue $ws.Range("D36") "4.59"
Set-TextValue $ws.Range("E36") "  -1.52%  "

Set-TextValue $ws.Range("E37") "  -0.59%  "

Set-TextValue $ws.Range("D38") "25.94"
Set-TextValue $ws.Range("E38") "  -5.36%  "

Set-TextValue $ws.Range("E39") "  -4.34%  "

Set-TextValue $ws.Range("E40") "  -3.59%  "

Set-TextValue $ws.Range("D41") "1.57"
Set-TextValue $ws.Range("E41") "  +7.49%  "

Set-TextValue $ws.Range("E42") "  -0.22%  "

Set-TextValue $ws.Range("D43") "0.680"
Set-TextValue $ws.Range("E43") "  +2.17%  "

Set-TextValue $ws.Range("D44") "3.120.50"
Set-TextValue $ws.Range("E44") "  -1.23%  "

Set-TextValue $ws.Range("D45") "36.85"
Set-TextValue $ws.Range("E45") "  -0.27%  "

Set-TextValue $ws.Range("D49") "0.989"
Set-TextValue $ws.Range("E49") "  +2.81%  "

Set-TextValue $ws.Range("D50") "6.06"
Set-TextValue $ws.Range("E50") "  +0.55%  "

Set-TextValue $ws.Range("D51") "20.46"
Set-TextValue $ws.Range("E51") "  -2.99%  "

# --- Row 47/48: Maker and VeChain swap positions with updated data ---
Set-TextValue $ws.Range("B47") "VeChain"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D47") "0.0261"
Set-TextValue $ws.Range("E47") "  +0.79%  "

Set-TextValue $ws.Range("B48") "Maker"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D48") "2.272.60"
Set-TextValue $ws.Range("E48") "  -0.99%  "
